# edit.ps1 - Apply the changes described in the commit:
# "feat: Hacer Ampliacion 2 Ejercicio 1 Jordan"
#
# Net effect inside the "Ejercicio 1" Ampliaciones list:
#  1. Paragraph "2. Permitir consultar el balance del contrato (0,5)"
#     changes its highlight from red to green.
#  2. A new list item is inserted right after it:
#     "3. Mostrar en la web el balance de BNB que hay en la wallet del
#     usuario (0,4)" (highlighted red, with a spell-check proof tag
#     around "wallet", matching the document's existing convention for
#     English technical terms).
#  3. The "Total:" score line changes from "0,4 puntos" to "0,9 puntos".

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Step 1: find the "2. Permitir consultar..." paragraph and the
# "Total: 0,4 puntos" paragraph by scanning Paragraphs once (text is
# stable / unique enough for -like matching).
# ---------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count

$idxPermitirConsultar = -1
$idxTotal = -1

for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    $txt = $p.Range.Text
    if ($idxPermitirConsultar -eq -1 -and $txt -like "*Permitir consultar el balance del contrato*") {
        $idxPermitirConsultar = $i
    }
    if ($idxTotal -eq -1 -and $txt -like "Total: 0,4*") {
        $idxTotal = $i
    }
}

if ($idxPermitirConsultar -eq -1) {
    throw "Could not find paragraph '2. Permitir consultar el balance del contrato (0,5)'"
}
if ($idxTotal -eq -1) {
    throw "Could not find paragraph 'Total: 0,4 puntos'"
}

# ---------------------------------------------------------------
# Step 2: recolor paragraph "2." from red to green (run + paragraph
# mark) by replacing its OOXML in place.
# ---------------------------------------------------------------
$p2 = $paras.Item($idxPermitirConsultar)
$xmlP2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="31"/></w:numPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>2. Permitir consultar el balance del contrato (0,5)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.Range.InsertXML($xmlP2) | Out-Null

# ---------------------------------------------------------------
# Step 3: insert the new "3. Mostrar en la web..." list item right
# after paragraph "2." (still red-highlighted, same list numbering).
# ---------------------------------------------------------------
$p2 = $paras.Item($idxPermitirConsultar)
$p2.Range.InsertParagraphAfter() | Out-Null
$newPara = $paras.Item($idxPermitirConsultar + 1)
$xmlNew3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="31"/></w:numPr><w:rPr><w:highlight w:val="red"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="red"/></w:rPr><w:t xml:space="preserve">3. Mostrar en la web el balance de BNB que hay en la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="red"/></w:rPr><w:t>wallet</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="red"/></w:rPr><w:t xml:space="preserve"> del usuario (0,4)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($xmlNew3) | Out-Null

# ---------------------------------------------------------------
# Step 4: update the "Total:" line from "0,4 puntos" to "0,9 puntos"
# (the paragraph index shifted down by one because of the insert
# above, so re-scan instead of reusing $idxTotal directly).
# ---------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count
$idxTotal2 = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    $txt = $p.Range.Text
    if ($txt -like "Total: 0,4*") {
        $idxTotal2 = $i
        break
    }
}
if ($idxTotal2 -eq -1) {
    throw "Could not find paragraph 'Total: 0,4 puntos' after insertion"
}
$pTotal = $paras.Item($idxTotal2)
$xmlTotal = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Total: </w:t></w:r><w:r><w:t>0</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t>9</w:t></w:r><w:r><w:t xml:space="preserve"> puntos</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pTotal.Range.InsertXML($xmlTotal) | Out-Null

Write-Output "Done."
